$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SAMPLE_TEST")

# New header cell L1 with the same text as the diff, then copy the
# header formatting (bold font, border, centered/top alignment) from
# the existing K1 header cell so it matches the rest of the header row.
$ws.Range("L1").Value = "time point"
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122)  # xlPasteFormats

# New data cells L2:L5, matching the existing "TP1" style (plain,
# unstyled inline-string cells like K2:K5).
$ws.Range("L2").Value = "TP1"
$ws.Range("L3").Value = "TP1"
$ws.Range("L4").Value = "TP1"
$ws.Range("L5").Value = "TP1"
